$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (A..T), per the Arf1-Chrm3 LR-pair update
$rowsData = @(
    @("ECs", "Arf1", "Chrm3", "ECs", 3, 1, 50.63667066666667, 151.910012, 0.20081482031288, 0.20081482031288, 3, 1, 3.259471333333333, 9.778414, 0.8750813810021594, 0.8750813810021593, 165.0487764534409, 1485.438988080968, 0.1757293102850956, 0.1757293102850955),
    @("ECs", "Arf1", "Chrm3", "FAPs", 3, 1, 50.63667066666667, 151.910012, 0.20081482031288, 0.20081482031288, 1, 0.3333333333333333, 0.1088336666666667, 0.326501, 0.02921894552415003, 0.02921894552415003, 5.510974536445778, 49.59877082801199, 0.005867597295164019, 0.005867597295164018),
    @("ECs", "Arf1", "Chrm3", "M2", 3, 1, 50.63667066666667, 151.910012, 0.20081482031288, 0.20081482031288, 2, 0.6666666666666666, 0.02652833333333333, 0.07958499999999999, 0.007122152090007321, 0.00712215209000732, 1.343306478335555, 12.08975830502, 0.001430233692195823, 0.001430233692195823),
    @("ECs", "Arf1", "Chrm3", "sCs", 3, 1, 50.63667066666667, 151.910012, 0.20081482031288, 0.20081482031288, 3, 1, 0.3299303333333334, 0.9897910000000001, 0.08857752138368333, 0.08857752138368331, 16.70657363194356, 150.359162687492, 0.01778767904042466, 0.01778767904042465),
    @("FAPs", "Arf1", "Chrm3", "ECs", 3, 1, 67.324, 201.972, 0.2669934019110801, 0.2669934019110801, 3, 1, 3.259471333333333, 9.778414, 0.8750813810021594, 0.8750813810021593, 219.4406480453333, 1974.965832408, 0.2336409548628126, 0.2336409548628125),
    @("FAPs", "Arf1", "Chrm3", "FAPs", 3, 1, 67.324, 201.972, 0.2669934019110801, 0.2669934019110801, 1, 0.3333333333333333, 0.1088336666666667, 0.326501, 0.02921894552415003, 0.02921894552415003, 7.327117774666666, 65.944059972, 0.007801265665747345, 0.007801265665747345),
    @("FAPs", "Arf1", "Chrm3", "M2", 3, 1, 67.324, 201.972, 0.2669934019110801, 0.2669934019110801, 2, 0.6666666666666666, 0.02652833333333333, 0.07958499999999999, 0.007122152090007321, 0.00712215209000732, 1.785993513333333, 16.07394162, 0.001901567615439164, 0.001901567615439164),
    @("FAPs", "Arf1", "Chrm3", "sCs", 3, 1, 67.324, 201.972, 0.2669934019110801, 0.2669934019110801, 3, 1, 0.3299303333333334, 0.9897910000000001, 0.08857752138368333, 0.08857752138368331, 22.21222976133334, 199.910067852, 0.02364961376708105, 0.02364961376708105),
    @("M2", "Arf1", "Chrm3", "ECs", 3, 1, 71.72398199999999, 215.171946, 0.2844428428612245, 0.2844428428612245, 3, 1, 3.259471333333333, 9.778414, 0.8750813810021594, 0.8750813810021593, 233.782263241516, 2104.040369173644, 0.2489106357471805, 0.2489106357471805),
    @("M2", "Arf1", "Chrm3", "FAPs", 3, 1, 71.72398199999999, 215.171946, 0.2844428428612245, 0.2844428428612245, 1, 0.3333333333333333, 0.1088336666666667, 0.326501, 0.02921894552415003, 0.02921894552415003, 7.805983948993998, 70.25385554094599, 0.008311119930296484, 0.008311119930296484),
    @("M2", "Arf1", "Chrm3", "M2", 3, 1, 71.72398199999999, 215.171946, 0.2844428428612245, 0.2844428428612245, 2, 0.6666666666666666, 0.02652833333333333, 0.07958499999999999, 0.007122152090007321, 0.00712215209000732, 1.90271770249, 17.12445932241, 0.002025845187771694, 0.002025845187771693),
    @("M2", "Arf1", "Chrm3", "sCs", 3, 1, 71.72398199999999, 215.171946, 0.2844428428612245, 0.2844428428612245, 3, 1, 0.3299303333333334, 0.9897910000000001, 0.08857752138368333, 0.08857752138368331, 23.663917289254, 212.975255603286, 0.02519524199597578, 0.02519524199597578),
    @("sCs", "Arf1", "Chrm3", "ECs", 3, 1, 62.47139133333334, 187.414174, 0.2477489349148154, 0.2477489349148154, 3, 1, 3.259471333333333, 9.778414, 0.8750813810021594, 0.8750813810021593, 203.6237092044485, 1832.613382840036, 0.2168004801070708, 0.2168004801070708),
    @("sCs", "Arf1", "Chrm3", "FAPs", 3, 1, 62.47139133333334, 187.414174, 0.2477489349148154, 0.2477489349148154, 1, 0.3333333333333333, 0.1088336666666667, 0.326501, 0.02921894552415003, 0.02921894552415003, 6.798990580574889, 61.190915225174, 0.007238962632942185, 0.007238962632942184),
    @("sCs", "Arf1", "Chrm3", "M2", 3, 1, 62.47139133333334, 187.414174, 0.2477489349148154, 0.2477489349148154, 2, 0.6666666666666666, 0.02652833333333333, 0.07958499999999999, 0.007122152090007321, 0.00712215209000732, 1.657261893087778, 14.91535703779, 0.001764505594600641, 0.00176450559460064),
    @("sCs", "Arf1", "Chrm3", "sCs", 3, 1, 62.47139133333334, 187.414174, 0.2477489349148154, 0.2477489349148154, 3, 1, 0.3299303333333334, 0.9897910000000001, 0.08857752138368333, 0.08857752138368331, 20.61120696640378, 185.500862697634, 0.02194498658020183, 0.02194498658020183)
)

$r = 2
foreach ($row in $rowsData) {
    for ($c = 0; $c -lt 20; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}
